$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Esoral" block (rows 2-12): Item Name (D) / UOM (E) pairs
$ws.Range("D2").Value = "Esoral Injection & Capsule 20"
$ws.Range("E2").Value = "Bundle"

$ws.Range("D3").Value = "Esoral 20mg Tablet"
$ws.Range("E3").Value = "20's"

$ws.Range("D4").Value = "Esoral 40mg Tablet"
$ws.Range("E4").Value = "20's"

$ws.Range("D5").Value = "Esoral Injection & MUPS 20"
$ws.Range("E5").Value = "Bundle"

$ws.Range("D6").Value = "Esoral 20mg Capsule 50's"
$ws.Range("E6").Value = "50's"

$ws.Range("D7").Value = "Esoral 20mg Tablet - 50's"
$ws.Range("E7").Value = "50's"

$ws.Range("D8").Value = "Esoral 20mg Tablet  80's"
$ws.Range("E8").Value = "80's"

$ws.Range("D9").Value = "Esoral 20mg  Tablet 100's"
$ws.Range("E9").Value = "100's"

$ws.Range("D10").Value = "Esoral 20mg Capsule 60's"
$ws.Range("E10").Value = "60's"

$ws.Range("D11").Value = "Esoral 20mg Capsule Container 30's"
$ws.Range("E11").Value = "30's"

$ws.Range("D12").Value = "Esoral 40mg EC Tablet - 42's"
$ws.Range("E12").Value = "42's"

# Reorder the "Losectil" block (rows 14-22): Item Name (D) / UOM (E) pairs
$ws.Range("D14").Value = "Losectil 20mg Capsule 500s"
$ws.Range("E14").Value = "500's"

$ws.Range("D17").Value = "Losectil 40mg Capsule (24's)"
$ws.Range("E17").Value = "24 's"

$ws.Range("D18").Value = "Losectil 20mg Powder for Oral Suspension"
$ws.Range("E18").Value = "20's"

$ws.Range("D20").Value = "Losectil 20mg Powder for Oral Suspension - 30's"
$ws.Range("E20").Value = "30's"

$ws.Range("D22").Value = "Losectil 20mg Capsule (100's)"
$ws.Range("E22").Value = "100 's"
